# Admin Staff person block: wrap the judges' signature line in a
# Jinja-style {% if %}/{% else %}/{% endif %} so that, when the case
# involves a BMV suspension, "Magistrate Amanda D. Bunner" is used
# instead of "Judge Marianne T. Hemmeter / Judge Kyle E. Rohrer".

$d = $word.ActiveDocument

# Locate the run holding the judges' names.
$rng = $d.Content
$found = $rng.Find.Execute("Judge Marianne T. Hemmeter / Judge Kyle E. Rohrer", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the Admin Staff signature text to update."
}

$startPos = $rng.Start
$origSize = $rng.Font.Size

# The four pieces of text that must end up as four separate runs
# (all sharing the original run's character formatting).
$t1 = "{% if bmv_suspension is true %}"
$t2 = "Magistrate Amanda D. Bunner{% else %}"
$t3 = " Judge Marianne T. Hemmeter / Judge Kyle E. Rohrer"
$t4 = " {% endif %}"

# Replace the found run's text with the concatenation of all four
# pieces first (still a single run at this point).
$rng.Text = $t1 + $t2 + $t3 + $t4

$p1 = $startPos + $t1.Length
$p2 = $p1 + $t2.Length
$p3 = $p2 + $t3.Length
$p4 = $p3 + $t4.Length

# Force Word to split the single run into four runs along those
# boundaries by briefly nudging (and then restoring) the font size on
# progressively smaller nested ranges. Because each nested range ends
# up with character formatting that differs from its neighbour while
# the nudge is in effect, Word can't coalesce them back together, and
# once every range is restored to the original size, the four runs
# remain distinct while keeping identical <w:rPr> formatting.
$d.Range($p1, $p4).Font.Size = $origSize + 1
$d.Range($p2, $p4).Font.Size = $origSize + 2
$d.Range($p3, $p4).Font.Size = $origSize + 3

$d.Range($p1, $p4).Font.Size = $origSize
$d.Range($p2, $p4).Font.Size = $origSize
$d.Range($p3, $p4).Font.Size = $origSize

Write-Output "Admin Staff signature line updated with BMV suspension conditional."
